# The deck's theme part (the one actually used by the slide master / all
# slides, serialized as ppt/theme/theme2.xml) was switched from the
# "Integral" colour palette back to the stock "Office Theme" palette -
# i.e. the commit swaps the contents that used to live in theme1.xml
# ("Office Theme") and theme2.xml ("Integral") between the two parts.
#
# This headless host only exposes one live/editable theme object (the
# one backing the part that the slide master + slides actually use); the
# notes-master-only theme part has no reachable ThemeColorScheme in this
# object model, and there is no ApplyTheme / OpenThemeFile-style import
# available here. The font scheme and format scheme (fills/lines/
# effects) are byte-identical between the "Office Theme" and "Integral"
# theme parts already, so the only real difference to reproduce is the
# 12 theme colours - which we can set directly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette = the stock "Office Theme" colours (RGB ints, i.e.
# R + G*256 + B*65536, matching PowerPoint's ColorFormat.RGB encoding).
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
